# Generate Report for Handback
# Adds the newly handed-back file 'Test`1.md' to the Overview and zh-cn
# report sheets, and refreshes the "Latest HO Xliff Generate Date" stamp
# on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = 'Test`1.md'
$wsOverview.Range("B3").Value = 'test\Test`1.md'
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = 'N\A'
$wsOverview.Range("G3").Value = "2017-11-06 03:16:53"
$wsOverview.Range("G3").NumberFormat = $dateFmt

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), 'https://github.com/OpenLocalizationTestOrg/LocaleLowerCaseTest/blob/2566724fb77af3b0efacc67e473603dec92ccf17/test/Test%601.md', "", "", 'test\Test`1.md')

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Refresh the handoff timestamp for the existing LocaleLowerCaseTest row.
$wsZhCn.Range("H2").Value = "2017-11-06 03:16:53"
$wsZhCn.Range("H2").NumberFormat = $dateFmt

# New row for the handed-back 'Test`1.md' file.
$wsZhCn.Range("A3").Value = 'Test`1.md'
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "test"
$wsZhCn.Range("E3").Value = "'"
$wsZhCn.Range("F3").Value = "'True"
$wsZhCn.Range("G3").Value = 'Test`1.2566724fb77af3b0efacc67e473603dec92ccf17.zh-cn.xlf'
$wsZhCn.Range("H3").Value = "2017-11-06 03:16:53"
$wsZhCn.Range("H3").NumberFormat = $dateFmt
$wsZhCn.Range("I3").Value = "'"
$wsZhCn.Range("J3").Value = 'Test`1.md'
$wsZhCn.Range("K3").Value = 'Test`1.2566724fb77af3b0efacc67e473603dec92ccf17.zh-cn.xlf'
$wsZhCn.Range("L3").Value = "2017-11-06 04:37:51"
$wsZhCn.Range("L3").NumberFormat = $dateFmt
$wsZhCn.Range("M3").Value = "'"
$wsZhCn.Range("N3").Value = "'"
$wsZhCn.Range("O3").Value = "'True"
$wsZhCn.Range("P3").Value = "'"
$wsZhCn.Range("Q3").Value = "'False"
$wsZhCn.Range("R3").Value = "'"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), 'https://github.com/OpenLocalizationTestOrg/LocaleLowerCaseTest/blob/2566724fb77af3b0efacc67e473603dec92ccf17/test/Test%601.md', "", "", 'Test`1.md')

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("J3"), 'https://github.com/OpenLocalizationTestOrg/LocaleLowerCaseTest.zh-cn/blob/2566724fb77af3b0efacc67e473603dec92ccf17/test/Test%601.md', "", "", 'Test`1.md')

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:R3"))

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Only the handoff timestamp is refreshed; no handback has happened yet
# for 'Test`1.md' on this locale, so no new row is added.
$wsDeDe.Range("H2").Value = "2017-11-06 03:16:56"
$wsDeDe.Range("H2").NumberFormat = $dateFmt
